$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24; this pushes the existing rows 24-47
# down to 25-48 and extends the used range to A1:R48.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new weekly price observation.
$ws.Range("A24").Value = 1
$ws.Range("B24").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C24").Value = "Arica y Parinacota"
$ws.Range("D24").Value = "2023-10-26"
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = 100112013
$ws.Range("G24").Value = "Alcachofa"
$ws.Range("H24").Value = "Madrigal"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 160
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 14500
$ws.Range("N24").Value = "`$/caja 40 unidades"
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 362
$ws.Range("Q24").Value = 40
$ws.Range("R24").Value = "Hortaliza"
